$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("clock")

# Insert a new row at 13 ("origin"), which pushes width-mount/offset-right/
# offset-left (and everything below) down by one row, and Excel will
# automatically re-point all formulas that referenced the old row numbers.
$ws.Rows("13:13").Insert()

# Populate the new row 13 with the "origin" label/value.
$ws.Range("A13").Value = "origin"
$ws.Range("B13").Value = 15.875

# Update offset-right (now row 15) and offset-left (now row 16) values.
$ws.Range("B15").Value = 6.08
$ws.Range("B16").Value = -5.92

# The bottom-most formula (now in row 24) needs to additionally add the new
# "origin" value (B13) on top of the auto-shifted B14-B16 expression.
$ws.Range("B24").Formula = "=B14-B16+B13"

# Leave the same cell selected/active as in the saved workbook.
$ws.Range("B23").Select() | Out-Null

$wb.Save()
